# NuAudio presentation update
# - Tweak wording on the "Backlog" slide
# - Insert a new "Fluid speech" slide before the "Dialogue example" slide
# - Revise some of the dialogue wording on the "Dialogue example" slide
# - Insert a new "Itunes Search API" slide after the "Dialogue example" slide
# - Refresh the cached date placeholders on the layouts/master

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Backlog slide: "iterator" -> "iterate"
# ---------------------------------------------------------------------------
$backlog = $p.Slides.Item(2)
$backlogBody = $backlog.Shapes.Item(2).TextFrame.TextRange
$backlogBody.Text = "Connecting itunes database to alexa responses`r" + `
    "Being able to iterate over more than just one response from alexa (MoreIntent)`r" + `
    "Add a user library, or repeated user request list, that allows alexa to search for song releases that the user might be interested in hearing (SuggestedIntent)`r" + `
    "`r`r"

# ---------------------------------------------------------------------------
# 2. Insert the new "Fluid speech" slide right before "Dialogue example"
#    (which is currently slide 4) -- it becomes the new slide 4.
# ---------------------------------------------------------------------------
$titleContentLayout = $p.SlideMaster.CustomLayouts.Item(2)
$fluidSlide = $p.Slides.AddSlide(4, $titleContentLayout)
$fluidSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Fluid speech"
$fluidBody = $fluidSlide.Shapes.Item(2).TextFrame.TextRange
$fluidBody.Text = "We are working on implementing SSML into our Alexa Responses`r" + `
    "Using SSML, we plan to make Alexa Responses more fluid`r" + `
    "Re-Writing some of our scripts, to make it seem more human like`r" + `
    "Adding a random response generator for alexa, so its not so automated, and is more natural responses "

# ---------------------------------------------------------------------------
# 3. "Dialogue example" slide is now slide 5 -- rework a few lines of
#    dialogue.
# ---------------------------------------------------------------------------
$dialogue = $p.Slides.Item(5)
$dialogueBody = $dialogue.Shapes.Item(2).TextFrame.TextRange

$q = [char]0x201C
$uq = [char]0x201D

$lines = @(
    "E-`t${q}Welcome to NuAudio! I can tell you about some new album releases:${uq}`t${q}Which artist would you like to know about?${uq}",
    "U-`t${q}who has albums coming out next week?${uq}",
    "E-`t${q}Rihanna has an album coming out November 7th.${uq}",
    "E-`t${q}Would you like to know about another artist?${uq}",
    "U-`t${q}Yes${uq}",
    "E-`t${q}Dope! What would you like to hear about now?${uq}",
    "U-`t${q}Help${uq}",
    "E-`t${q}I can tell you about artists or I can tell you about new music releases by date. What would you like to hear?${uq}",
    "U-`t${q}Tell me about Slip knot${uq}",
    "E-`t${q}Slip knot released This Old Dog on May 5, 2017${uq}`t${q}Would you like to know about another?${uq}",
    "U-`t${q}No.${uq}",
    "E-`t${q}Thanks for using NuAudio. Goodbye!${uq}"
)
$dialogueBody.Text = [string]::Join("`r", $lines)

# Restore the superscript "th" in "November 7th."
$plain = $dialogueBody.Text
$thIdx = $plain.IndexOf("7th.")
if ($thIdx -ge 0) {
    $thChars = $dialogueBody.Characters($thIdx + 2, 2)
    $thChars.Font.Superscript = $true
}

# ---------------------------------------------------------------------------
# 4. Insert the new "Itunes Search API" slide right after "Dialogue example"
#    (which is slide 5) -- it becomes the new slide 6.
# ---------------------------------------------------------------------------
$itunesSlide = $p.Slides.AddSlide(6, $titleContentLayout)
$itunesSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Itunes Search API"
$itunesBody = $itunesSlide.Shapes.Item(2).TextFrame.TextRange
$itunesBody.Text = "Successfully able to send an http get request to itunes api.`r" + `
    "Able to get a response from api in the form of a Json object.`r" + `
    "Able to parse Json object for an artists 5 most recent songs and latest album`r" + `
    ""

# ---------------------------------------------------------------------------
# 5. Refresh the cached "today" date shown on every layout + the master
#    (PowerPoint recalculates these automatically; we just bump the cached
#    text forward a day to mirror the authoring session).
# ---------------------------------------------------------------------------
try {
    $master = $p.SlideMaster
    for ($i = 1; $i -le $master.Shapes.Placeholders.Count; $i++) {
        $ph = $master.Shapes.Placeholders.Item($i)
        if ($ph.PlaceholderFormat.Type -eq 16) {
            if ($ph.TextFrame.TextRange.Text -eq "11/8/17") {
                $ph.TextFrame.TextRange.Text = "11/9/17"
            }
        }
    }
} catch {
}

for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    try {
        $layout = $p.SlideMaster.CustomLayouts.Item($i)
        for ($j = 1; $j -le $layout.Shapes.Placeholders.Count; $j++) {
            $ph = $layout.Shapes.Placeholders.Item($j)
            if ($ph.PlaceholderFormat.Type -eq 16) {
                if ($ph.TextFrame.TextRange.Text -eq "11/8/17") {
                    $ph.TextFrame.TextRange.Text = "11/9/17"
                }
            }
        }
    } catch {
    }
}

# ---------------------------------------------------------------------------
# 6. Touch the notes master so it gets materialised alongside the deck
#    (mirrors enabling Notes Page view in the authoring session).
# ---------------------------------------------------------------------------
try {
    $nm = $p.NotesMaster
} catch {
}
